# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1 (13:22 -> 13:52)
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 13:52"

# Province/city labels that moved to a different row because the data
# got re-sorted by total cases (column B) descending.
$ws.Range("A19").Value = "Valladolid"
$ws.Range("A20").Value = "Salamanca"
$ws.Range("A25").Value = "Segovia"
$ws.Range("A26").Value = "Cantabria"
$ws.Range("A27").Value = "Granada"
$ws.Range("A29").Value = "Leon"
$ws.Range("A30").Value = "Pontevedra"
$ws.Range("A39").Value = "Soria"
$ws.Range("A40").Value = "Aragon"
$ws.Range("A41").Value = "Avila"
$ws.Range("A42").Value = "Cuenca"

# Updated case counts (Casos totales, Casos activos, Recuperados, Muertes)
$ws.Range("B19").Value = 2123
$ws.Range("C19").Value = 804
$ws.Range("D19").Value = 1124
$ws.Range("E19").Value = 195

$ws.Range("B20").Value = 2115
$ws.Range("C20").Value = 579
$ws.Range("D20").Value = 1292
$ws.Range("E20").Value = 244

$ws.Range("B25").Value = 1760
$ws.Range("C25").Value = 490
$ws.Range("D25").Value = 1131
$ws.Range("E25").Value = 139

$ws.Range("B26").Value = 1719
$ws.Range("C26").Value = 281
$ws.Range("D26").Value = 1331
$ws.Range("E26").Value = 107

$ws.Range("B27").Value = 1686
$ws.Range("C27").Value = 251
$ws.Range("D27").Value = 1287
$ws.Range("E27").Value = 148

$ws.Range("B29").Value = 1538
$ws.Range("C29").Value = 751
$ws.Range("D29").Value = 551
$ws.Range("E29").Value = 236

$ws.Range("B30").Value = 1536
$ws.Range("C30").Value = 333
$ws.Range("D30").Value = 1411
$ws.Range("E30").Value = 30

$ws.Range("B33").Value = 1170
$ws.Range("C33").Value = 492
$ws.Range("D33").Value = 542
$ws.Range("E33").Value = 136

$ws.Range("B39").Value = 927
$ws.Range("C39").Value = 228
$ws.Range("D39").Value = 619
$ws.Range("E39").Value = 80

$ws.Range("B40").Value = 907
$ws.Range("C40").Value = 29
$ws.Range("D40").Value = 838
$ws.Range("E40").Value = 40

$ws.Range("B41").Value = 897
$ws.Range("C41").Value = 347
$ws.Range("D41").Value = 456
$ws.Range("E41").Value = 94

$ws.Range("B42").Value = 874
$ws.Range("C42").Value = 2205
$ws.Range("D42").Value = 9768
$ws.Range("E42").Value = 129

$ws.Range("B45").Value = 605
$ws.Range("C45").Value = 171
$ws.Range("D45").Value = 388
$ws.Range("E45").Value = 46

$ws.Range("B50").Value = 408
$ws.Range("C50").Value = 148
$ws.Range("D50").Value = 209
$ws.Range("E50").Value = 51
